$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 299-300 (existing rows 299+ shift down to 301+)
$ws.Rows("299:300").Insert()

# New row 299 data
$d299 = Get-Date -Year 2021 -Month 11 -Day 16 -Hour 0 -Minute 0 -Second 0
$ws.Range("A299").Value = 4
$ws.Range("B299").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C299").Value = "Los Lagos"
$ws.Range("D299").Value = $d299
$ws.Range("E299").Value = 10
$ws.Range("F299").Value = "Fruta"
$ws.Range("G299").Value = 100102
$ws.Range("H299").Value = "Cítricos"
$ws.Range("I299").Value = 100102003
$ws.Range("J299").Value = "Limón"
$ws.Range("K299").Value = "Sin especificar"
$ws.Range("L299").Value = "1a amarillo"
$ws.Range("M299").Value = 1200
$ws.Range("N299").Value = 12000
$ws.Range("O299").Value = 12500
$ws.Range("P299").Value = 12250
$ws.Range("Q299").Value = "$/malla 18 kilos"
$ws.Range("R299").Value = "Provincia de Melipilla"
$ws.Range("S299").Value = 681
$ws.Range("T299").Value = 18

# New row 300 data
$d300 = Get-Date -Year 2021 -Month 11 -Day 16 -Hour 0 -Minute 0 -Second 0
$ws.Range("A300").Value = 4
$ws.Range("B300").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C300").Value = "Los Lagos"
$ws.Range("D300").Value = $d300
$ws.Range("E300").Value = 10
$ws.Range("F300").Value = "Fruta"
$ws.Range("G300").Value = 100102
$ws.Range("H300").Value = "Cítricos"
$ws.Range("I300").Value = 100102003
$ws.Range("J300").Value = "Limón"
$ws.Range("K300").Value = "Sin especificar"
$ws.Range("L300").Value = "2a amarillo"
$ws.Range("M300").Value = 400
$ws.Range("N300").Value = 10500
$ws.Range("O300").Value = 10500
$ws.Range("P300").Value = 10500
$ws.Range("Q300").Value = "$/malla 18 kilos"
$ws.Range("R300").Value = "Provincia de Melipilla"
$ws.Range("S300").Value = 583
$ws.Range("T300").Value = 18
